$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Telephone number and email address of the applicant."
$ws.Range("B25").Value = "Name and contact information for the parties making the application."
$ws.Range("B31").Value = "Name and contact information if an agent is being used."
$ws.Range("B35").Value = "Name and contact information if an agent is being used."
$ws.Range("B43").Value = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
$ws.Range("B44").Value = "Signed and dated verification of the application's accuracy."
$ws.Range("B47").Value = "Details of any hedgerows being removed as part of the development"
$ws.Range("B53").Value = "Details of pre-application advice received from the planning authority"
$ws.Range("B58").Value = "Where the proposed development will be built."
$ws.Range("B67").Value = "Information to help the planning authority arrange a site visit"
